# Actualización desde MV -datos-
# Adds 7 new daily rows (08-09-2021 .. 16-09-2021) to the bottom of the
# "Diaria" sheet, continuing the existing Serie/Cupo/Monto... table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row data (row -> Serie date, Cupo, Monto demandado, Total adjudicado,
#     Monto adjudicado bancos, Monto adjudicado AFP, Tasa) ---------------------
$rows = @(
    @{ R=79; Date="08-09-2021"; B=100000; C=$null;  D=0;      E=$null;  F=$null;  G=$null },
    @{ R=80; Date="09-09-2021"; B=100000; C=275000; D=100000; E=100000; F=0;      G=2     },
    @{ R=81; Date="10-09-2021"; B=100000; C=$null;  D=0;      E=$null;  F=$null;  G=$null },
    @{ R=82; Date="13-09-2021"; B=50000;  C=140000; D=50000;  E=45000;  F=5000;   G=2     },
    @{ R=83; Date="14-09-2021"; B=50000;  C=50000;  D=50000;  E=40000;  F=10000;  G=2     },
    @{ R=84; Date="15-09-2021"; B=50000;  C=80000;  D=50000;  E=40000;  F=10000;  G=2.03  },
    @{ R=85; Date="16-09-2021"; B=50000;  C=$null;  D=0;      E=$null;  F=$null;  G=$null }
)

# Column A holds the "Serie" date as plain text (matches the existing column
# of "dd-mm-yyyy" text labels). Assigning it through .Value/.Formula directly
# lets Excel's smart-parser reinterpret short day/month pairs (day <= 12) as
# real dates, so we instead enter it as a formula producing the literal
# string and immediately convert that formula to its static value with a
# copy / paste-values round-trip -- the same result a user gets by typing
# ="08-09-2021" and then Paste Special > Values, with no left-over
# number-format / style changes.
foreach ($item in $rows) {
    $cell = $ws.Cells.Item($item.R, 1)
    $cell.Formula = '="' + $item.Date + '"'
}

$dateRange = $ws.Range("A79:A85")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

foreach ($item in $rows) {
    $r = $item.R
    $ws.Cells.Item($r, 2).Value = $item.B
    if ($item.C -ne $null) { $ws.Cells.Item($r, 3).Value = $item.C }
    $ws.Cells.Item($r, 4).Value = $item.D
    if ($item.E -ne $null) { $ws.Cells.Item($r, 5).Value = $item.E }
    if ($item.F -ne $null) { $ws.Cells.Item($r, 6).Value = $item.F }
    if ($item.G -ne $null) { $ws.Cells.Item($r, 7).Value = $item.G }
}
